# Auto-generated Excel COM-interop script
# Applies the 2024-08-03 violent crime data update across the workbook:
# updates the 2024 (column K) year-to-date totals (and a couple of 2023
# column J reclassifications) on the Citywide Totals, By Neighborhood,
# and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{}

$sheetUpdates['Citywide Totals'] = @{
    'K2' = 4667
    'K3' = 4793
    'J4' = 1830
    'K4' = 986
    'K6' = 5396
    'J7' = 29300
    'K7' = 16185
}

$sheetUpdates['By Neighborhood'] = @{
    'K5' = 37
    'K7' = 478
    'K8' = 1091
    'K12' = 29
    'K15' = 164
    'K20' = 375
    'K23' = 167
    'K24' = 50
    'K25' = 79
    'K27' = 147
    'K29' = 860
    'K31' = 182
    'K33' = 683
    'K37' = 546
    'K42' = 603
    'K47' = 103
    'K48' = 203
    'K50' = 83
    'K51' = 205
    'K52' = 423
    'K53' = 215
    'K54' = 311
    'K58' = 11
    'K63' = 50
    'K64' = 102
    'K65' = 370
    'K69' = 36
    'K73' = 139
    'K75' = 55
    'K78' = 191
    'K79' = 399
    'K80' = 57
    'K83' = 343
    'K84' = 118
    'K85' = 739
    'J86' = 177
    'K86' = 107
    'K88' = 187
    'K91' = 172
    'K95' = 283
    'K98' = 80
    'K99' = 272
    'J101' = 29300
    'K101' = 16185
}

$sheetUpdates['Auburn Gresham'] = @{
    'K3' = 154
    'K7' = 478
}

$sheetUpdates['South Shore'] = @{
    'K3' = 247
    'K4' = 43
    'K6' = 170
    'K7' = 739
}

$sheetUpdates['Little Village'] = @{
    'K6' = 159
    'K7' = 423
}

$sheetUpdates['Norwood Park'] = @{
    'K2' = 10
    'K7' = 36
}

$sheetUpdates['Logan Square'] = @{
    'K2' = 54
    'K7' = 215
}

$sheetUpdates['Austin'] = @{
    'K2' = 303
    'K3' = 328
    'K4' = 64
    'K6' = 367
    'K7' = 1091
}

$sheetUpdates['South Chicago'] = @{
    'K3' = 123
    'K7' = 343
}

$sheetUpdates['Garfield Park'] = @{
    'K2' = 189
    'K3' = 255
    'K6' = 195
    'K7' = 683
}

$sheetUpdates['West Pullman'] = @{
    'K3' = 96
    'K7' = 283
}

$sheetUpdates['Grand Crossing'] = @{
    'K2' = 154
    'K3' = 182
    'K7' = 546
}

$sheetUpdates['New City'] = @{
    'K6' = 145
    'K7' = 370
}

$sheetUpdates['Woodlawn'] = @{
    'K3' = 112
    'K4' = 17
    'K7' = 272
}

$sheetUpdates['Gage Park'] = @{
    'K2' = 59
    'K7' = 182
}

$sheetUpdates['South Deering'] = @{
    'K6' = 26
    'K7' = 118
}

$sheetUpdates['Loop'] = @{
    'K6' = 161
    'K7' = 311
}

$sheetUpdates['Englewood'] = @{
    'K2' = 245
    'K3' = 307
    'K7' = 860
}

$sheetUpdates['Lake View'] = @{
    'K4' = 29
    'K7' = 203
}

$sheetUpdates['Humboldt Park'] = @{
    'K2' = 163
    'K3' = 188
    'K7' = 603
}

$sheetUpdates['Rogers Park'] = @{
    'K4' = 20
    'K7' = 191
}

$sheetUpdates['Dunning'] = @{
    'K2' = 16
    'K6' = 16
    'K7' = 50
}

$sheetUpdates['Douglas'] = @{
    'K6' = 44
    'K7' = 167
}

$sheetUpdates['Washington Park'] = @{
    'K2' = 42
    'K7' = 172
}

$sheetUpdates['Roseland'] = @{
    'K6' = 98
    'K7' = 399
}

$sheetUpdates['Near South Side'] = @{
    'K2' = 20
    'K4' = 12
    'K6' = 38
    'K7' = 102
}

$sheetUpdates['Chicago Lawn'] = @{
    'K2' = 126
    'K3' = 119
    'K4' = 16
    'K7' = 375
}

$sheetUpdates['East Side'] = @{
    'K2' = 25
    'K7' = 79
}

$sheetUpdates['Kenwood'] = @{
    'K3' = 31
    'K7' = 103
}

$sheetUpdates['Brighton Park'] = @{
    'K6' = 50
    'K7' = 164
}

$sheetUpdates['Wicker Park'] = @{
    'K3' = 16
    'K7' = 80
}

$sheetUpdates['Lincoln Square'] = @{
    'K2' = 21
    'K7' = 83
}

$sheetUpdates['Portage Park'] = @{
    'K3' = 35
    'K7' = 139
}

$sheetUpdates['United Center'] = @{
    'K2' = 46
    'K6' = 77
    'K7' = 187
}

$sheetUpdates['Armour Square'] = @{
    'K6' = 16
    'K7' = 37
}

$sheetUpdates['Edgewater'] = @{
    'K2' = 38
    'K4' = 17
    'K7' = 147
}

$sheetUpdates['Streeterville'] = @{
    'J4' = 96
    'K4' = 44
    'J7' = 177
    'K7' = 107
}

$sheetUpdates['Pullman'] = @{
    'K3' = 18
    'K7' = 55
}

$sheetUpdates['Little Italy, UIC'] = @{
    'K6' = 70
    'K7' = 205
}

$sheetUpdates['Rush & Division'] = @{
    'K3' = 14
    'K6' = 27
    'K7' = 57
}

$sheetUpdates['Beverly'] = @{
    'K4' = 3
    'K7' = 29
}

$sheetUpdates['Millenium Park'] = @{
    'K4' = 1
    'K6' = 9
    'K7' = 11
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetUpdates[$sheetName]
    foreach ($ref in $cellMap.Keys) {
        $ws.Range($ref).Value = $cellMap[$ref]
    }
}
